# Working time log: fill in the missing Nov-15 entry (row 9) and refresh
# the sheet's navigation/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 was present (with a live shared formula in D9) but A9:C9 were
# --- still blank. Add the missing Date / Start Time / End Time values.

# Date column - reuse the existing date formatting (row 3) via copy/paste
# of formats so the new cell picks up the same style as the rest of the
# table instead of minting a redundant number format.
$ws.Range("A9").Value = 43054
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Start / End time columns - plain time-of-day values with the default
# time display.
$ws.Range("B9").Value = 0.33333333333333331
$ws.Range("C9").Value = 0.36458333333333331
$ws.Range("B9:C9").NumberFormat = "h:mm"

# D9's shared formula (=ABS(C9-B9)) and the Table1 totals row recalculate
# automatically once the inputs above are in place.

# --- Refresh the view: scroll the window back to the top of the sheet and
# --- move the selection off the old out-of-range cell.
$ws.Range("B10").Select() | Out-Null
